# Update the "想去人数" (want-to-go count) figures on the 展览 and 全部类型
# sheets, reflecting the refreshed scrape output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 6916
$wsExpo.Range("F5").Value = 50
$wsExpo.Range("F6").Value = 1072
$wsExpo.Range("F7").Value = 161
$wsExpo.Range("F8").Value = 6

# --- Sheet "全部类型" ------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6916
$wsAll.Range("F5").Value = 50
$wsAll.Range("F6").Value = 1072
$wsAll.Range("F7").Value = 161
$wsAll.Range("F9").Value = 6
